# Generate Report for Handoff
#
# The localization CI pipeline re-ran and produced a new report: the
# status moves from "In Translation" to "Ready for handoff" and the
# "Latest HO Xliff Generate Date" / "Latest Handoff Datetime" timestamps
# are refreshed. The Status/date columns are widened slightly so the new
# values are fully visible.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# Width (in "characters") that renders to the widened column size used by
# the refreshed report (matches the new, wider Status/date columns).
$newColWidth = 16.3333333333333

# --- Overview sheet ------------------------------------------------------
# E2/F2 hold the zh-cn / de-de status, G2 the latest handoff-xliff-generate
# timestamp.
$overview.Range("E2").Value = "Ready for handoff"
$overview.Range("F2").Value = "Ready for handoff"
$overview.Range("G2").Value = "2016-08-31 01:02:34"

$overview.Columns.Item(5).ColumnWidth = $newColWidth
$overview.Columns.Item(6).ColumnWidth = $newColWidth

# --- zh-cn sheet -----------------------------------------------------------
$zhcn.Range("C2").Value = "Ready for handoff"
$zhcn.Range("H2").Value = "2016-08-31 01:02:29"

$zhcn.Columns.Item(3).ColumnWidth = $newColWidth

# --- de-de sheet -----------------------------------------------------------
$dede.Range("C2").Value = "Ready for handoff"
$dede.Range("H2").Value = "2016-08-31 01:02:34"

$dede.Columns.Item(3).ColumnWidth = $newColWidth
